$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "Version" column between Company (B) and Number of Rating (old C, now D)
$ws.Columns("C").Insert()
$ws.Range("C1").Value = "Version"

# Fill in Version values for the existing apps (force text so "4.7.16" isn't
# auto-parsed as a date by the smart-entry logic)
$ws.Range("C2").Value = "2.2.0"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "4.7.16"
$ws.Range("C3").Style = "Normal"
$ws.Range("C4").Value = "2.1.1"

# Fix the Mi Home rating (was 2, now 1.9)
$ws.Range("E3").Value = 1.9

# Add two new app rows
$ws.Range("A5").Value = "Cosmote"
$ws.Range("B5").Value = "COSMOTE"
$ws.Range("C5").Value = "1.3.0"
$ws.Range("D5").Value = "N/A"
$ws.Range("E5").Value = 4
$ws.Range("E5").NumberFormat = "0.0"

$ws.Range("A6").Value = "ImperiHome"
$ws.Range("B6").Value = "Evertygo"
$ws.Range("C6").Value = "4.1.2"
$ws.Range("D6").Value = 106
$ws.Range("E6").Value = 4.5
$ws.Range("E6").NumberFormat = "0.0"

# Adjust column widths: App Name (A) widens to fit "ImperiHome", new Version
# column (C) matches the Company column's width, Number of Rating (D) keeps
# its original width automatically after the column insert.
$ws.Columns("A").ColumnWidth = 12.25
$ws.Columns("C").ColumnWidth = 21.92

# Restore the selected/active cell as it was left after editing
$ws.Range("D11").Select() | Out-Null
